$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Table1")
$newCol = $tbl.ListColumns.Add(3)
$newCol.Name = "Trạng thái"
Write-Host $tbl.ListColumns.Count
Write-Host $tbl.Range.Address()
